# [PHOENIX-5910] changes in search trade license and approval details
#
# The "grievanceDetails" sheet (Sheet2) had a typo fixed in the grievance
# description text, and the sheet's active view/selection was updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grievanceDetails")
$ws.Activate()

# Fix spelling: "Deadly deseases ..." -> "Deadly diseases ..."
$ws.Range("D2").Value = "Deadly diseases are spreading because of mosquitoes"

# Update the sheet view: scroll so column E is the left-most visible column,
# and move the active selection to H19.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("H19").Select()
